$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the UUID-like ID values in column A for rows 2-4
$ws.Range("A2").Value = "80e2d59c-96d4-449e-b326-51f9ce876cac"
$ws.Range("A3").Value = "f2e0af68-d87b-4ad2-8ca9-ea027f0f0498"
$ws.Range("A4").Value = "a96a2d42-4301-4541-9b3a-442cbe5dfc25"

# Update the Price for row 3 (Dropdown B) from 15 to 50
$ws.Range("D3").Value = 50

# Delete row 5 (Widget E) entirely
$ws.Rows(5).Delete()
